$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1255
$ws.Range("F4").Value = 57
$ws.Range("F5").Value = 5547
$ws.Range("F6").Value = 1777
$ws.Range("F7").Value = 6346
$ws.Range("G7").Value = 80
$ws.Range("F8").Value = 139
$ws.Range("F9").Value = 1910
$ws.Range("F10").Value = 512
$ws.Range("F11").Value = 4
$ws.Range("F13").Value = 29
$ws.Range("G15").Value = 45
$ws.Range("I15").Value = '//i0.hdslb.com/bfs/openplatform/202409/ms9IIHAn1725447474436.jpeg'
$ws.Range("F17").Value = 7881
$ws.Range("F18").Value = 7881
$ws.Range("F19").Value = 139
$ws.Range("F21").Value = 180
$ws.Range("F28").Value = 46
$ws.Range("F29").Value = 172
$ws.Range("F30").Value = 1733
$ws.Range("F31").Value = 799
$ws.Range("F32").Value = 368
$ws.Range("F35").Value = 78
$ws.Range("F36").Value = 89
$ws.Range("F37").Value = 3915

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 360
$ws.Range("F20").Value = 33

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2272
$ws.Range("F4").Value = 682
$ws.Range("F5").Value = 265

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2272
$ws.Range("F4").Value = 682
$ws.Range("F5").Value = 1255
$ws.Range("F7").Value = 57
$ws.Range("F9").Value = 360
$ws.Range("F10").Value = 5547
$ws.Range("F11").Value = 265
$ws.Range("F12").Value = 1777
$ws.Range("F13").Value = 6346
$ws.Range("G13").Value = 80
$ws.Range("F14").Value = 139
$ws.Range("F15").Value = 1910
$ws.Range("F17").Value = 512
$ws.Range("F19").Value = 29
$ws.Range("C21").Value = '桐庐·唯泽动漫游戏嘉年华'
$ws.Range("D21").Value = '城南路277号 桐庐海博大酒店'
$ws.Range("E21").Value = '2024.09.22 10:00-09.22 17:00'
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 45
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=91601'
$ws.Range("I21").Value = '//i0.hdslb.com/bfs/openplatform/202409/ms9IIHAn1725447474436.jpeg'
$ws.Range("F23").Value = 7881
$ws.Range("F24").Value = 7881
$ws.Range("F25").Value = 139
$ws.Range("F27").Value = 180
$ws.Range("F33").Value = 46
$ws.Range("F34").Value = 172
$ws.Range("F35").Value = 1733
$ws.Range("F36").Value = 799
$ws.Range("F38").Value = 368
$ws.Range("F45").Value = 3915
$ws.Range("F46").Value = 33
